$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text "Active Cases" -> "Active cases"
$ws.Range("B1").Value = "Active cases"

$names = @(
    "3398 BlueCross Elly Kay Mordialloc",
    "3601 Baptcare Westhaven community",
    "3647 Aurrum Aged Care Reservoir",
    "3653 Fronditha Thalpori St Albans Aged Care",
    "3975 Aurrum Aged Care Brunswick West",
    "4257 BlueCross The Gables Camberwell",
    "4295 Hope Aged Care Sunshine West",
    "4314 Estia Health Ardeer",
    "44095 Myrniong Primary School Myrniong",
    "44404 Castlemaine North Primary School Castlemaine",
    "44490 Armadale Primary School Armadale",
    "44507 Inverloch Primary School Inverloch",
    "44593 Torquay P-6 College Torquay",
    "44642 Irymple South Primary School Irymple South",
    "44745 Briar Hill Primary School Briar Hill",
    "44753 Brunswick North West Primary School Brunswick West",
    "4479 Whittlesea Lodge Whittlesea",
    "44799 Eastwood Primary School Ringwood East",
    "44855 Beverley Hills Primary School Doncaster East",
    "44893 Greenhills Primary School Greensborough",
    "44960 Thomastown West Primary School",
    "45013 Gladstone Views Primary School",
    "45147 Maramba Primary School Narre Warren",
    "45168 Ranfurly Primary School Mildura",
    "45257 Roxburgh Rise Primary School Roxburgh Park",
    "45305 Lockington Consolidated SchoolLockington",
    "4574 Village Glen Aged Care Residences Mornington",
    "45757 Saint Joseph's Primary School Warragul",
    "45858 St Bernard's Primary Coburg",
    "45901 St Mary's Parish Primary School",
    "45958 Ave Maria College Aberfeldie Workplace",
    "46074 St Justin's Catholic Primary School Wheelers Hill",
    "46078 Corpus Christi Primary School Werribee",
    "46086 St Kevin's Primary School Hampton Park",
    "46320 St Mary's Coptic Orthodox College Coolaroo",
    "46327 Victory Christian College Strathdale",
    "50279 Dallas Brooks Community Primary School Dallas",
    "50412 Geelong Grammar School Bostock House Campus Newtown",
    "52390 Our Lady of the Way Catholic Primary School Wallan",
    "52694 Pakenham Primary School Pakenham",
    "52777 Mirripoa Primary School Mount Duneed School Camp",
    "Cardinia Waters Retirement Village Pakenham",
    "Confirmed Omicron Sircuit Bar Fitzroy",
    "Confirmed Omicron Variant The Peel HotelCollingwood",
    "Flight QF10 LHR to MEL via DRW 13 Dec",
    "Green Gables Lodge Warburton",
    "Greendale Hotel Greendale",
    "JBS Australia Brooklyn",
    "Jackson's Hotel Age of Love Event Toorak",
    "Kororoit Creek Primary School Burnside Heights Oct-Dec",
    "PGL Camp Rumbug Foster North",
    "Social Gathering 11 Dec Windsor",
    "Social Gathering 11 December Fitzroy",
    "St Brigid's Parish Primary School Mordialloc",
    "St Clare's Primary School Officer",
    "St Vincents Hospital Melbourne Emergency Department Fitzroy",
    "StarTrack- Melbourne Tullamarine",
    "The George Lounge St Kilda",
    "Thomastown West Primary School Camp Doxa's Malmsbury"
)

$counts = @(28,13,11,20,13,16,14,17,19,63,13,10,25,14,16,31,15,29,11,12,18,30,10,36,15,42,10,14,10,10,19,11,14,14,10,12,10,11,49,21,21,16,16,20,13,27,24,40,13,20,36,23,10,14,11,11,22,15,16)

for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
    $ws.Cells.Item($i + 2, 2).Value = $counts[$i]
}

# Remove now-unused trailing rows (61-64)
$ws.Range("A61:B64").EntireRow.Delete()

